$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new weekly data row above the current row 113, pushing the
# existing rows 113-117 down to 114-118.
$ws.Rows("113:113").Insert()

# Populate the new row 113 with the latest weekly quote.
$ws.Range("A113").Value = 4
$ws.Range("B113").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C113").Value = "Los Lagos"
$ws.Range("D113").Value = 44939
$ws.Range("E113").Value = 10
$ws.Range("F113").Value = 100112031
$ws.Range("G113").Value = "Poroto verde"
$ws.Range("H113").Value = "Magnum"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 40
$ws.Range("K113").Value = 35000
$ws.Range("L113").Value = 35000
$ws.Range("M113").Value = 35000
$ws.Range("N113").Value = "$/malla 25 kilos"
$ws.Range("O113").Value = "Región Metropolitana"
$ws.Range("P113").Value = 1400
$ws.Range("Q113").Value = 25
$ws.Range("R113").Value = "Hortaliza"
